$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (avoids Excel
# auto-converting numeric-looking strings like "709.72" into numbers),
# then clear the temporary number-format so the cell keeps its original
# (default) style, matching the source data which stores these as plain text.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "71.110.43"
$ws.Range("E2").Value = "  +2.97%  "

Set-TextValue $ws.Range("D3") "3.804.71"
$ws.Range("E3").Value = "  +0.82%  "

$ws.Range("E4").Value = "  +0.07%  "

Set-TextValue $ws.Range("D5") "709.72"
$ws.Range("E5").Value = "  +12.69%  "

Set-TextValue $ws.Range("D6") "173.01"
$ws.Range("E6").Value = "  +4.49%  "

Set-TextValue $ws.Range("D7") "3.804.19"
$ws.Range("E7").Value = "  +0.88%  "

$ws.Range("E8").Value = "  +0.00%  "

Set-TextValue $ws.Range("D9") "0.529"
$ws.Range("E9").Value = "  +1.20%  "

Set-TextValue $ws.Range("D10") "0.164"
$ws.Range("E10").Value = "  +3.42%  "

Set-TextValue $ws.Range("D11") "7.51"
$ws.Range("E11").Value = "  +10.10%  "

Set-TextValue $ws.Range("D12") "0.462"
$ws.Range("E12").Value = "  +1.22%  "

Set-TextValue $ws.Range("D13") "0.0000262"
$ws.Range("E13").Value = "  +9.55%  "

Set-TextValue $ws.Range("D14") "36.31"
$ws.Range("E14").Value = "  +4.10%  "

Set-TextValue $ws.Range("D15") "4.445.30"
$ws.Range("E15").Value = "  +0.82%  "

Set-TextValue $ws.Range("D16") "3.807.39"
$ws.Range("E16").Value = "  +0.62%  "

Set-TextValue $ws.Range("D17") "71.122.49"
$ws.Range("E17").Value = "  +2.96%  "

Set-TextValue $ws.Range("D18") "17.93"
$ws.Range("E18").Value = "  +1.35%  "

Set-TextValue $ws.Range("D19") "7.24"
$ws.Range("E19").Value = "  +3.16%  "

$ws.Range("E20").Value = "  +0.54%  "

$ws.Range("E21").Value = "  +18.05%  "

Set-TextValue $ws.Range("D22") "485.29"
$ws.Range("E22").Value = "  +3.74%  "

Set-TextValue $ws.Range("D23") "0.716"
$ws.Range("E23").Value = "  +1.82%  "

$ws.Range("E24").Value = "  +3.47%  "

Set-TextValue $ws.Range("D25") "83.89"
$ws.Range("E25").Value = "  +2.20%  "

Set-TextValue $ws.Range("D26") "12.44"
$ws.Range("E26").Value = "  +2.65%  "

Set-TextValue $ws.Range("D27") "10.57"
$ws.Range("E27").Value = "  +4.24%  "

$ws.Range("E28").Value = "  +2.80%  "

Set-TextValue $ws.Range("D29") "3.955.42"
$ws.Range("E29").Value = "  +0.81%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D30") "3.16"
$ws.Range("E30").Value = "  +18.46%  "

$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D31") "0.999"
$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D32") "7.59"
$ws.Range("E32").Value = "  +6.74%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D33") "2.30"
$ws.Range("E33").Value = "  +0.55%  "

Set-TextValue $ws.Range("D34") "29.67"
$ws.Range("E34").Value = "  +4.38%  "

Set-TextValue $ws.Range("D35") "0.180"
$ws.Range("E35").Value = "  +0.59%  "

Set-TextValue $ws.Range("D36") "9.26"
$ws.Range("E36").Value = "  +4.00%  "

Set-TextValue $ws.Range("D37") "0.998"
$ws.Range("E37").Value = "  -0.09%  "

Set-TextValue $ws.Range("D38") "3.754.72"
$ws.Range("E38").Value = "  +0.77%  "

$ws.Range("E39").Value = "  +2.02%  "

Set-TextValue $ws.Range("D40") "3.52"
$ws.Range("E40").Value = "  +7.68%  "

Set-TextValue $ws.Range("D41") "5.99"
$ws.Range("E41").Value = "  +3.33%  "

Set-TextValue $ws.Range("D42") "2.23"
$ws.Range("E42").Value = "  +11.45%  "

Set-TextValue $ws.Range("D43") "0.000331"
$ws.Range("E43").Value = "  +26.07%  "

$ws.Range("E44").Value = "  +0.76%  "

Set-TextValue $ws.Range("D45") "0.999"
$ws.Range("E45").Value = "  -0.12%  "

$ws.Range("E46").Value = "  -0.04%  "

Set-TextValue $ws.Range("D47") "162.22"
$ws.Range("E47").Value = "  +3.87%  "

Set-TextValue $ws.Range("D48") "49.50"
$ws.Range("E48").Value = "  +5.30%  "

Set-TextValue $ws.Range("D49") "45.10"
$ws.Range("E49").Value = "  +2.67%  "

$ws.Range("E50").Value = "  -1.09%  "

$ws.Range("E51").Value = "  +2.64%  "

Write-Host "Updated cryptos list"